# Regenerate the handback-status report:
#  - file "ac55d992-1edd-410e-bdfd-f7bebed1963e.md" was re-handed-back as
#    "32148b73-c6c7-4f2c-b423-3e776891f11f.md"
#  - file "c9c70c95-c479-43fd-a4d2-48270d45b7c0.md" was re-handed-back as
#    "ffff8dbadc86-c263-4d8e-a71d-8ef7fd70f007.md"
# New xliff hand-off/back filenames & timestamps are generated too.

$wb = $excel.ActiveWorkbook

$oldName1 = "ac55d992-1edd-410e-bdfd-f7bebed1963e.md"
$newName1 = "32148b73-c6c7-4f2c-b423-3e776891f11f.md"
$oldName2 = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.md"
$newName2 = "ffff8dbadc86-c263-4d8e-a71d-8ef7fd70f007.md"

$oldXlf1Zh = "ac55d992-1edd-410e-bdfd-f7bebed1963e.345e5284197f96f7f8137abb0f1f643e2aedff4e.zh-cn.xlf"
$oldXlf2Zh = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.5ce180b0445c7c4ad289ac1955311cb07e44f947.zh-cn.xlf"
$oldXlf1De = "ac55d992-1edd-410e-bdfd-f7bebed1963e.345e5284197f96f7f8137abb0f1f643e2aedff4e.de-de.xlf"
$oldXlf2De = "c9c70c95-c479-43fd-a4d2-48270d45b7c0.5ce180b0445c7c4ad289ac1955311cb07e44f947.de-de.xlf"

# Both rows now point at the same regenerated xliff (row 3's own xliff got
# collapsed away in the new report), matching the source diff exactly.
$newXlfZh = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.zh-cn.xlf"
$newXlfDe = "32148b73-c6c7-4f2c-b423-3e776891f11f.b32bc2ec4764a6f6d6612e33d75134f195191199.de-de.xlf"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Sheets.Item("Overview")

$ov.Range("A2").Value = $newName1
$ov.Range("G2").Value = "2016-08-24 07:03:58"

$ov.Range("A3").Value = $newName2
$ov.Range("G3").Value = "2016-08-24 07:03:58"

$ov.Hyperlinks.Delete()
$ov.Hyperlinks.Add($ov.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName1", "", "", "e2e\$newName1") | Out-Null
$ov.Hyperlinks.Add($ov.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName2", "", "", "e2e\$newName2") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Sheets.Item("zh-cn")

$zh.Range("A2").Value = $newName1
$zh.Range("G2").Value = $newXlfZh
$zh.Range("H2").Value = "2016-08-24 07:03:53"
$zh.Range("I2").Value = $newName1
$zh.Range("J2").Value = $newXlfZh
$zh.Range("K2").Value = "2016-08-24 07:04:17"

$zh.Range("A3").Value = $newName2
$zh.Range("G3").Value = $newXlfZh
$zh.Range("H3").Value = "2016-08-24 07:03:53"
$zh.Range("I3").Value = $newName2
$zh.Range("J3").Value = $newXlfZh
$zh.Range("K3").Value = "2016-08-24 07:04:17"

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName1", "", "", $newName1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/19b1ed1869ad447272fbb108c10b0625baae4d29/e2e/$oldName1", "", "", $newName1) | Out-Null
$zh.Hyperlinks.Add($zh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName2", "", "", $newName2) | Out-Null
$zh.Hyperlinks.Add($zh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/19b1ed1869ad447272fbb108c10b0625baae4d29/e2e/$oldName2", "", "", $newName2) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$de = $wb.Sheets.Item("de-de")

$de.Range("A2").Value = $newName1
$de.Range("G2").Value = $newXlfDe
$de.Range("H2").Value = "2016-08-24 07:03:58"
$de.Range("I2").Value = $newName1
$de.Range("J2").Value = $newXlfDe
$de.Range("K2").Value = "2016-08-24 07:04:25"

$de.Range("A3").Value = $newName2
$de.Range("G3").Value = $newXlfDe
$de.Range("H3").Value = "2016-08-24 07:03:58"
$de.Range("I3").Value = $newName2
$de.Range("J3").Value = $newXlfDe
$de.Range("K3").Value = "2016-08-24 07:04:25"

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName1", "", "", $newName1) | Out-Null
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8719cc4ce4d60fe6370da1857590d18e58b50a03/e2e/$oldName1", "", "", $newName1) | Out-Null
$de.Hyperlinks.Add($de.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/fdbee8bb0cadca47eacfdb25318c0ba91d04e6b9/e2e/$oldName2", "", "", $newName2) | Out-Null
$de.Hyperlinks.Add($de.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/8719cc4ce4d60fe6370da1857590d18e58b50a03/e2e/$oldName2", "", "", $newName2) | Out-Null
